$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.703.88"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.554.53"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("E6").Value = "  +6.95%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +9.44%  "
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").Value = "2.547.68"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "42.769.16"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +8.06%  "
$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +8.08%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +1.93%  "
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("E36").Value = "  +5.52%  "
$ws.Range("E37").Value = "  +15.42%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("B40").Value = "ApeXProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E40").Value = "  +30.25%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "2.058.79"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("E47").Value = "  +5.46%  "
$ws.Range("E48").Value = "  +10.84%  "
$ws.Range("D49").Value = "2.805.37"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E51").Value = "  +2.47%  "

$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "302.21"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "97.51"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "0.544"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "36.02"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "7.51"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "13.48"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "6.59"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "71.56"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "255.87"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "28.06"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "38.88"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "6.00"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "155.78"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "0.0803"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "3.32"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "25.90"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "18.36"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "2.06"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "3.85"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "0.0304"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "0.999"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "88.04"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "76.07"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Value = "103.57"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163, -4142, $false, $false)
$scratch.Clear()
$excel.CutCopyMode = $false
